$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New row of data (row 6), mirroring the format of row 5
$rowData = @(42861, 229, 309, 67, 3555, 3, 202, 216, 657, 130, 1, 225, 10, 0, 10, 83, 15, 1)

for ($col = 1; $col -le $rowData.Length; $col++) {
    $ws.Cells.Item(6, $col).Value = $rowData[$col - 1]
}

# Match the date style used in column A for rows 2-5 (numFmtId 14 => short date)
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A6").Value = 42861

# Update selection to match the new active cell / selection range
$ws.Range("A6:R6").Select()
